# Fruta / hortaliza, semanal
#
# The underlying data rows (2-20) get reshuffled as part of the weekly
# refresh: each destination row ends up holding the full set of column
# values (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Precio $/Kg, Kg/unidad, etc.) that used to live in a different row.
# Row 1 (headers) is untouched, and every row keeps its own formatting -
# only the underlying values move.
#
# $newRow -> $oldRow : row $newRow should end up with the values that
# used to be in row $oldRow.
$rowMap = @{
    2  = 6
    3  = 7
    4  = 16
    5  = 11
    6  = 9
    7  = 10
    8  = 20
    9  = 18
    10 = 3
    11 = 4
    12 = 13
    13 = 14
    14 = 2
    15 = 5
    16 = 12
    17 = 19
    18 = 8
    19 = 17
    20 = 15
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 20
$firstCol = 1
$lastCol = 20

# 1) Snapshot every source row's values BEFORE any writes happen, so the
#    permutation doesn't clobber data we still need to read.
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $original[$r] = $rowVals
}

# 2) Write each destination row using the snapshotted source row's values.
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $rowVals = $original[$oldRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $rowVals[$c - 1]
    }
}
